# Update the header date.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-09-15 Friday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2023-09-16 Saturday", 2) | Out-Null

# Update the multiplication answers in the table. Cells are addressed
# directly by (row, column) and their Range.Text is replaced in place so
# that identical-looking cells (e.g. the two "68×67=4556" cells) are each
# updated independently and correctly, instead of relying on a
# document-wide Find/Replace which would touch every matching occurrence
# at once.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "60×49=2940"
$t.Cell(1, 2).Range.Text = "80×95=7600"
$t.Cell(1, 3).Range.Text = "49×55=2695"
$t.Cell(1, 4).Range.Text = "64×32=2048"
$t.Cell(1, 5).Range.Text = "99×12=1188"

$t.Cell(5, 1).Range.Text = "50×82=4100"
$t.Cell(5, 2).Range.Text = "69×96=6624"
$t.Cell(5, 3).Range.Text = "66×18=1188"
$t.Cell(5, 4).Range.Text = "90×83=7470"
$t.Cell(5, 5).Range.Text = "82×12=984"

$t.Cell(10, 1).Range.Text = "41×15=615"
$t.Cell(10, 2).Range.Text = "24×14=336"
$t.Cell(10, 3).Range.Text = "74×46=3404"
$t.Cell(10, 4).Range.Text = "14×80=1120"
$t.Cell(10, 5).Range.Text = "71×85=6035"

$t.Cell(15, 1).Range.Text = "42×67=2814"
$t.Cell(15, 2).Range.Text = "34×54=1836"
$t.Cell(15, 3).Range.Text = "13×19=247"
$t.Cell(15, 4).Range.Text = "44×56=2464"
$t.Cell(15, 5).Range.Text = "80×45=3600"

$t.Cell(20, 1).Range.Text = "70×35=2450"
$t.Cell(20, 2).Range.Text = "74×37=2738"
$t.Cell(20, 3).Range.Text = "22×25=550"
$t.Cell(20, 4).Range.Text = "53×84=4452"
$t.Cell(20, 5).Range.Text = "40×61=2440"

Write-Output "edit complete"
